# Implements "Implemented Aug 7 2023 feedback" - converts the chained
# {% endif %}{% if ... %} Jinja blocks into proper {% elif ... %} blocks,
# and drops a stray "and" in the Property branch.

$d = $word.ActiveDocument

$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”
$rsq = [char]0x2019  # ’

# 1) "...{% endif %}{% if poa_type == “Property” %}Durable Power of Attorney
#     for Property and dated..."
#    -> "...{% elif poa_type == “Property” %}Durable Power of Attorney for
#        Property dated..."
$old1 = "{% endif %}{% if poa_type == " + $lq + "Property" + $rq + " %}Durable Power of Attorney for Property and dated"
$new1 = "{% elif poa_type == " + $lq + "Property" + $rq + " %}Durable Power of Attorney for Property dated"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) "...{% endif %}{% if poa_type == “Both” %}Durable..."
#    -> "...{% elif poa_type == “Both” %}Durable..."
$old2 = "{% endif %}{% if poa_type == " + $lq + "Both" + $rq + " %}Durable"
$new2 = "{% elif poa_type == " + $lq + "Both" + $rq + " %}Durable"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) "...{{health_end_date}}.{% endif %}{%if poa_type == “Property” %}
#     {{property_end_date}}.{% endif %}{% if poa_type == “Both” %}for..."
#    -> "...{{health_end_date}}.{%elif poa_type == “Property” %}
#        {{property_end_date}}. {% elif poa_type == “Both” %}for..."
$old3 = "{{health_end_date}}.{% endif %}{%if poa_type == " + $lq + "Property" + $rq + " %}{{property_end_date}}.{% endif %}{% if poa_type == " + $lq + "Both" + $rq + " %}for"
$new3 = "{{health_end_date}}.{%elif poa_type == " + $lq + "Property" + $rq + " %}{{property_end_date}}. {% elif poa_type == " + $lq + "Both" + $rq + " %}for"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# The "_GoBack" bookmark marks the author's last edit point. The edits above
# touched the "Property" branch text last, so re-anchor "_GoBack" right after
# "...Durable Power of Attorney for Property " in that branch (it previously
# sat after the same phrase in the "Both" branch).
$anchor = "{% elif poa_type == " + $lq + "Property" + $rq + " %}Durable Power of Attorney for Property "
$r = $d.Content
$r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

Write-Output $d.Content.Text
